$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 2017 <-> 2021 OWASP Top 10 mapping text (ES localisation fixes) ---
$ws.Range("C7").Value  = "A01:2017-Injection"
$ws.Range("E7").Value  = "A01:2021-Pérdida de Control de Acceso"

$ws.Range("C8").Value  = "A02:2017-Broken Authentication"
$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"

$ws.Range("C9").Value  = "A03:2017-Sensitive Data Exposure"
$ws.Range("E9").Value  = "A03:2021-Inyección"

$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("D10").Value = "(Nueva)"
$ws.Range("E10").Value = "A04:2021-Diseño Inseguro"

$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"

$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"

$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("E13").Value = "A07:2021-Fallas de Identificación y Autenticación"

$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("E14").Value = "A08:2021-Fallas en el Software y en la Integridad de los Datos"

$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"

$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"
$ws.Range("D16").Value = "(Nueva)"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"

$ws.Range("E17").Value = "* A partir de la encuesta"

# --- Column C narrowed now that the English 2017 labels are shorter ---
$ws.Columns.Item(3).ColumnWidth = 47.6666666666667

# --- Remove the two stray/duplicate connector shapes left over in the drawing ---
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()

# --- Update the saved selection / scroll position ---
$ws.Range("E23").Select()
